# "damage on successful save"
#
# Inserts a new "saved_damage" column just before the existing
# "attack_damage" column (column S) on both the Heroes and Monsters
# sheets. The old "attack_damage" column's data is shifted one column
# to the right (S -> T), and the new column S is populated with a
# header ("saved_damage") plus a numeric value on the monster/boss row
# (row 4) representing the damage still dealt on a successful save.

$wb = $excel.ActiveWorkbook

# --- Heroes sheet ---
$heroes = $wb.Worksheets.Item("Heroes")
$heroes.Activate()

# Shift the existing "attack_damage" column (S) one column to the right.
$heroes.Range("T1").Value = $heroes.Range("S1").Value2
$heroes.Range("T2").Value = $heroes.Range("S2").Value2
$heroes.Range("T3").Value = $heroes.Range("S3").Value2
$heroes.Range("T4").Value = $heroes.Range("S4").Value2

# Column S becomes the new "saved_damage" column.
$heroes.Range("S1").Value = "saved_damage"
$heroes.Range("S2").ClearContents()
$heroes.Range("S3").ClearContents()
$heroes.Range("S4").Value = 0.5

$heroes.Range("S1:T4").Select()

# --- Monsters sheet ---
$monsters = $wb.Worksheets.Item("Monsters")
$monsters.Activate()

# Shift the existing "attack_damage" column (S) one column to the right.
$monsters.Range("T1").Value = $monsters.Range("S1").Value2
$monsters.Range("T2").Value = $monsters.Range("S2").Value2
$monsters.Range("T3").Value = $monsters.Range("S3").Value2
$monsters.Range("T4").Value = $monsters.Range("S4").Value2

# Column S becomes the new "saved_damage" column.
$monsters.Range("S1").Value = "saved_damage"
$monsters.Range("S2").ClearContents()
$monsters.Range("S3").ClearContents()
$monsters.Range("S4").Value = 0

$monsters.Range("S5").Select()
